$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1499.6
$ws.Range("J19").Value = 1624.5
$ws.Range("L19").Value = 1624.5
$ws.Range("N19").Value = -1974.5
$ws.Range("H41").Value = 1450.2222
$ws.Range("I41").Value = 1231.1666
$ws.Range("J41").Value = 1888.3334
$ws.Range("K41").Value = 1231.1666
$ws.Range("L41").Value = 1888.3334
$ws.Range("M41").Value = -791.1666
$ws.Range("N41").Value = -2768.3334
$ws.Range("H103").Value = 941.2105
$ws.Range("I103").Value = 1115.4445
$ws.Range("J103").Value = 784.4
$ws.Range("K103").Value = 3346.3335
$ws.Range("L103").Value = 2353.2
$ws.Range("M103").Value = -2760.3335
$ws.Range("N103").Value = -3525.2
$ws.Range("H125").Value = 3197.3076
$ws.Range("J125").Value = 3830.75
$ws.Range("L125").Value = 34476.75
$ws.Range("N125").Value = -39396.75
$ws.Range("H132").Value = 2387.1177
$ws.Range("I132").Value = 2283.8572
$ws.Range("J132").Value = 2869
$ws.Range("K132").Value = 6851.571599999999
$ws.Range("L132").Value = 8607
$ws.Range("M132").Value = -4321.571599999999
$ws.Range("N132").Value = -13667
$ws.Range("H138").Value = 1161978.5
$ws.Range("J138").Value = 1548584.9
$ws.Range("L138").Value = 4645754.699999999
$ws.Range("N138").Value = -4656034.699999999
$ws.Range("H141").Value = 3235.353
$ws.Range("I141").Value = 2008.6154
$ws.Range("K141").Value = 6025.8462
$ws.Range("M141").Value = -845.8462

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5187.1875
$ws.Range("I61").Value = 3375.875
$ws.Range("K61").Value = 3375.875
$ws.Range("M61").Value = -3163.875
$ws.Range("H63").Value = 4146.231
$ws.Range("I63").Value = 3090.1
$ws.Range("K63").Value = 3090.1
$ws.Range("M63").Value = -2404.1
$ws.Range("H66").Value = 4146.231
$ws.Range("I66").Value = 3090.1
$ws.Range("K66").Value = 15450.5
$ws.Range("M66").Value = -12018.5
$ws.Range("H74").Value = 1438.4286
$ws.Range("I74").Value = 1139.6
$ws.Range("K74").Value = 1139.6
$ws.Range("M74").Value = -265.5999999999999
$ws.Range("H77").Value = 1438.4286
$ws.Range("I77").Value = 1139.6
$ws.Range("K77").Value = 5698
$ws.Range("M77").Value = -1330
$ws.Range("H101").Value = 53995.5
$ws.Range("J101").Value = 53995.5
$ws.Range("L101").Value = 53995.5
$ws.Range("N101").Value = -60485.5
$ws.Range("H114").Value = 39999
$ws.Range("J114").Value = 39999
$ws.Range("L114").Value = 39999
$ws.Range("N114").Value = -48677
$ws.Range("H124").Value = 67500
$ws.Range("J124").Value = 67500
$ws.Range("L124").Value = 67500
$ws.Range("N124").Value = -77320
$ws.Range("H125").Value = 67500
$ws.Range("J125").Value = 67500
$ws.Range("L125").Value = 67500
$ws.Range("N125").Value = -77340
$ws.Range("H136").Value = 5187.1875
$ws.Range("I136").Value = 3375.875
$ws.Range("K136").Value = 10127.625
$ws.Range("M136").Value = -7577.625

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 965.55554
$ws.Range("I64").Value = 999
$ws.Range("K64").Value = 999
$ws.Range("M64").Value = -774
$ws.Range("H67").Value = 965.55554
$ws.Range("I67").Value = 999
$ws.Range("K67").Value = 999
$ws.Range("M67").Value = -219

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4096.75
$ws.Range("I16").Value = 1500
$ws.Range("J16").Value = 4962.3335
$ws.Range("K16").Value = 1500
$ws.Range("L16").Value = 4962.3335
$ws.Range("M16").Value = -1213
$ws.Range("N16").Value = -5536.3335
$ws.Range("H58").Value = 3313.2144
$ws.Range("I58").Value = 3106.6155
$ws.Range("J58").Value = 5999
$ws.Range("K58").Value = 3106.6155
$ws.Range("L58").Value = 5999
$ws.Range("M58").Value = -2903.6155
$ws.Range("N58").Value = -6405
$ws.Range("H68").Value = 20000
$ws.Range("I68").Value = 20000
$ws.Range("J68").Value = 20000
$ws.Range("K68").Value = 20000
$ws.Range("L68").Value = 20000
$ws.Range("M68").Value = -19251
$ws.Range("N68").Value = -21498
$ws.Range("H71").Value = 20000
$ws.Range("I71").Value = 20000
$ws.Range("J71").Value = 20000
$ws.Range("K71").Value = 60000
$ws.Range("L71").Value = 60000
$ws.Range("M71").Value = -56256
$ws.Range("N71").Value = -67488
$ws.Range("H94").Value = 2209.5217
$ws.Range("I94").Value = 2667.8572
$ws.Range("K94").Value = 2667.8572
$ws.Range("M94").Value = -2216.8572
$ws.Range("H99").Value = 20771.428
$ws.Range("I99").Value = 20399.7
$ws.Range("J99").Value = 21109.363
$ws.Range("K99").Value = 20399.7
$ws.Range("L99").Value = 21109.363
$ws.Range("M99").Value = -18901.7
$ws.Range("N99").Value = -24105.363
$ws.Range("H113").Value = 4096.75
$ws.Range("I113").Value = 1500
$ws.Range("J113").Value = 4962.3335
$ws.Range("K113").Value = 1500
$ws.Range("L113").Value = 4962.3335
$ws.Range("M113").Value = 670
$ws.Range("N113").Value = -9302.333500000001
$ws.Range("H122").Value = 3877.2666
$ws.Range("I122").Value = 3025.9
$ws.Range("J122").Value = 5580
$ws.Range("K122").Value = 9077.700000000001
$ws.Range("L122").Value = 16740
$ws.Range("M122").Value = -6627.700000000001
$ws.Range("N122").Value = -21640
$ws.Range("H126").Value = 20771.428
$ws.Range("I126").Value = 20399.7
$ws.Range("J126").Value = 21109.363
$ws.Range("K126").Value = 61199.10000000001
$ws.Range("L126").Value = 63328.08900000001
$ws.Range("M126").Value = -58729.10000000001
$ws.Range("N126").Value = -68268.08900000001
$ws.Range("H132").Value = 4244.143
$ws.Range("I132").Value = 3817.2856
$ws.Range("J132").Value = 4671
$ws.Range("K132").Value = 11451.8568
$ws.Range("L132").Value = 14013
$ws.Range("M132").Value = -8921.856800000001
$ws.Range("N132").Value = -19073
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").ClearContents()
$ws.Range("H136").Value = 3313.2144
$ws.Range("I136").Value = 3106.6155
$ws.Range("J136").Value = 5999
$ws.Range("K136").Value = 9319.8465
$ws.Range("L136").Value = 17997
$ws.Range("M136").Value = -6769.8465
$ws.Range("N136").Value = -23097

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 355.86667
$ws.Range("I38").Value = 35.333332
$ws.Range("K38").Value = 105.999996
$ws.Range("M38").Value = 241.000004
$ws.Range("H43").Value = 1666.6666
$ws.Range("I43").Value = 500
$ws.Range("J43").Value = 4000
$ws.Range("K43").Value = 1500
$ws.Range("L43").Value = 12000
$ws.Range("M43").Value = -1386
$ws.Range("N43").Value = -12228
$ws.Range("H75").Value = 532.3333
$ws.Range("J75").Value = 999
$ws.Range("L75").Value = 2997
$ws.Range("N75").Value = -4993
$ws.Range("H78").Value = 532.3333
$ws.Range("J78").Value = 999
$ws.Range("L78").Value = 8991
$ws.Range("N78").Value = -18975

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 44996
$ws.Range("J93").Value = 44996
$ws.Range("L93").Value = 44996
$ws.Range("N93").Value = -48740
$ws.Range("H108").Value = 119999
$ws.Range("J108").Value = 119999
$ws.Range("L108").Value = 119999
$ws.Range("N108").Value = -127679
$ws.Range("H113").Value = 3110.375
$ws.Range("I113").Value = 3110.375
$ws.Range("K113").Value = 3110.375
$ws.Range("M113").Value = -940.375
$ws.Range("H122").Value = 1497.6364
$ws.Range("I122").Value = 1509.375
$ws.Range("J122").Value = 1466.3334
$ws.Range("K122").Value = 4528.125
$ws.Range("L122").Value = 4399.0002
$ws.Range("M122").Value = -2078.125
$ws.Range("N122").Value = -9299.0002
$ws.Range("H123").Value = 50040.855
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 50040.855
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 50040.855
$ws.Range("M123").ClearContents()
$ws.Range("N123").Value = -54940.855
$ws.Range("H132").Value = 4767.5454
$ws.Range("I132").Value = 3053.5
$ws.Range("J132").Value = 5747
$ws.Range("K132").Value = 9160.5
$ws.Range("L132").Value = 17241
$ws.Range("M132").Value = -6630.5
$ws.Range("N132").Value = -22301

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 5120535.5
$ws.Range("I43").Value = 2106250
$ws.Range("J43").Value = 7381250
$ws.Range("K43").Value = 2106250
$ws.Range("L43").Value = 7381250
$ws.Range("M43").Value = -2106057
$ws.Range("N43").Value = -7381636
$ws.Range("H136").Value = 4232.032
$ws.Range("I136").Value = 3009.1304
$ws.Range("J136").Value = 7747.875
$ws.Range("K136").Value = 9027.3912
$ws.Range("L136").Value = 23243.625
$ws.Range("M136").Value = -6477.3912
$ws.Range("N136").Value = -28343.625

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 6328.8945
$ws.Range("J122").Value = 8496.5
$ws.Range("L122").Value = 25489.5
$ws.Range("N122").Value = -30389.5
$ws.Range("H126").Value = 3306.9285
$ws.Range("I126").Value = 2299.7273
$ws.Range("J126").Value = 7000
$ws.Range("K126").Value = 6899.1819
$ws.Range("L126").Value = 21000
$ws.Range("M126").Value = -4429.1819
$ws.Range("N126").Value = -25940
$ws.Range("H132").Value = 13961.223
$ws.Range("I132").Value = 10612
$ws.Range("J132").Value = 16640.6
$ws.Range("K132").Value = 31836
$ws.Range("L132").Value = 49921.8
$ws.Range("M132").Value = -29306
$ws.Range("N132").Value = -54981.8
$ws.Range("H136").Value = 5809.0356
$ws.Range("I136").Value = 5881.684
$ws.Range("J136").Value = 5655.6665
$ws.Range("K136").Value = 17645.052
$ws.Range("L136").Value = 16966.9995
$ws.Range("M136").Value = -15095.052
$ws.Range("N136").Value = -22066.9995
